# Change the "SHIPPING" line-item label to "Subtotal" in the invoice
# totals table (single occurrence in the document).
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "SHIPPING",  # FindText
    $true,       # MatchCase
    $true,       # MatchWholeWord
    $false,      # MatchWildcards
    $false,      # MatchSoundsLike
    $false,      # MatchAllWordForms
    $true,       # Forward
    1,           # Wrap (wdFindContinue)
    $false,      # Format
    "Subtotal",  # ReplaceWith
    2            # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find 'SHIPPING' text to replace with 'Subtotal'."
}

Write-Output "Replaced SHIPPING -> Subtotal: $found"
